$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17, shifting existing rows 17-42 down to 18-43.
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with the new data record.
$ws.Cells.Item(17, 1).Value = 10
$ws.Cells.Item(17, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(17, 3).Value = "La Araucanía"
$ws.Cells.Item(17, 4).Value = 44483
$ws.Cells.Item(17, 5).Value = 9
$ws.Cells.Item(17, 6).Value = 300000001
$ws.Cells.Item(17, 7).Value = "Rabanito"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7600
$ws.Cells.Item(17, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(17, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(17, 16).Value = 633
$ws.Cells.Item(17, 17).Value = 12
$ws.Cells.Item(17, 18).Value = "Hortaliza"

# Apply the same date-number-format style used by the other "Fecha" column cells.
$ws.Cells.Item(17, 4).NumberFormat = $ws.Cells.Item(18, 4).NumberFormat
